# Update crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.760.45"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "2.175.27"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.90"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.22%  "

$ws.Range("E14").Value = "  -3.04%  "

$ws.Range("D15").Value = "2.501.57"
$ws.Range("E15").Value = "  -2.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "2.174.40"
$ws.Range("E17").Value = "  -3.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.785"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.64%  "

$ws.Range("D19").Value = "41.635.42"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.00%  "

$ws.Range("E28").Value = "  -10.01%  "

$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0778"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.19%  "

$ws.Range("E36").Value = "  -3.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "

$ws.Range("E38").Value = "  -6.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.190"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.11%  "

$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.70%  "

$ws.Range("E48").Value = "  -4.03%  "

$ws.Range("E49").Value = "  -4.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.55%  "

$ws.Range("E51").Value = "  -2.17%  "
